$wb = $excel.ActiveWorkbook

# --- CategoricalVariables sheet -------------------------------------------
# Insert a new "below detection limit" (code 6) row for C_quality_flag,
# mirroring the existing one for N_quality_flag. This shifts every
# N_quality_flag row down by one.
$ws2 = $wb.Worksheets.Item("CategoricalVariables")
$ws2.Activate()

$ws2.Rows.Item(6).Insert()
$ws2.Range("A6").Value = "C_quality_flag"
$ws2.Range("B6").Value = 6
$ws2.Range("C6").Value = "below detection limit"
$ws2.Range("A7").Select()

# --- ColumnHeaders sheet ----------------------------------------------------
# Update the C_quality_flag attributeDefinition to mention the lab flag for
# values below detection (same wording already used for N_quality_flag).
$ws1 = $wb.Worksheets.Item("ColumnHeaders")
$ws1.Activate()

$ws1.Range("B15").Value = "Carbon sample IODE Quality Flag primary level for carbon data with lab flag for value below detection"
$ws1.Rows.Item(15).RowHeight = 28.8
$ws1.Range("B15").Select()
